# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
# Both sheets mirror the same source rows, so the same F-column updates
# apply to each.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    2  = 4915
    5  = 807
    6  = 248
    10 = 217
    15 = 4343
    16 = 6624
    22 = 4070
    23 = 434
    25 = 41
    26 = 2658
    28 = 543
    31 = 340
    36 = 1003
    39 = 73
    43 = 85
    44 = 614
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
